$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Ramesh"
$ws.Range("B2").Value = "Kumar"
$ws.Range("C2").Value = "ramesh.kumar@yopmail.com"
$ws.Range("D2").Value = 9876543210

$ws.Range("A3").Value = "Gita"
$ws.Range("B3").Value = "Ben"
$ws.Range("C3").Value = "gita.ben@yopmail.com"
$ws.Range("D3").Value = 9012345678

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:ramesh.kumar@yopmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:gita.ben@yopmail.com")

$ws.Columns("C").ColumnWidth = 26.6
$ws.Range("E3").Select() | Out-Null
